# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" sheets:
#   F4: 44  -> 45
#   F5: 270 -> 271

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 45
    $ws.Range("F5").Value = 271
}
